# Applies the "Lc 9" edit:
#   1. Paragraph 1 (title "Lc 9,18-22..."): drop the spell-check proofErr
#      markers around "Lc", merge "Lc" + " 9,18-22" into a single red
#      (EE0000) bold run, colour the following line-break run the same
#      red, and split "UNA vez que..." into "Erase una" + " vez que...".
#   2. Paragraph 5 ("Quien es Jesus para ti?"): drop proofErr markers,
#      merge into a single run.
#   3. Paragraph 21 ("Esta pregunta tambien..."): drop proofErr markers,
#      merge the runs that were only split apart because of them.
#
# Word's Find/Replace keeps stray <w:proofErr> markers behind when a
# run that sits between a spellStart/spellEnd pair is merged away, so
# we rebuild the three affected paragraphs from scratch with
# Range.InsertXML, which lets us hand Word the exact run layout we
# want (including the new <w:color/> runs) in one shot.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $bodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------
# Paragraph 1 : "Lc 9,18-22" title block
# ---------------------------------------------------------------
$para1 = '<w:p w14:paraId="65EC5A3B" w14:textId="43C79DE2" w:rsidR="001F22DF" w:rsidRDefault="001F22DF" w:rsidP="001F22DF">' +
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/><w:color w:val="EE0000"/></w:rPr><w:t>Lc 9,18-22</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/><w:color w:val="EE0000"/></w:rPr><w:br/></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Erase una</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> vez que Jes&#250;s estaba orando solo, lo acompa&#241;aban sus disc&#237;pulos y les pregunt&#243;:</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#171;&#191;Qui&#233;n dice la gente que soy yo?&#187;.</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>Ellos contestaron:</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#171;Unos, que Juan el Bautista; otros, que El&#237;as, otros dicen que ha resucitado uno de los antiguos profetas&#187;.</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#201;l les pregunt&#243;:</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#171;Y ustedes, &#191;qui&#233;n dicen que soy yo?&#187;.</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>Pedro respondi&#243;:</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#171;El Mes&#237;as de Dios&#187;.</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#201;l les prohibi&#243; terminantemente dec&#237;rselo a nadie. Porque dec&#237;a:</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>&#171;El Hijo del hombre tiene que padecer mucho, ser desechado por los ancianos, sumos sacerdotes y escribas, ser ejecutado y resucitar al tercer d&#237;a&#187;.</w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:br/><w:t>Palabra del Se&#241;or.</w:t></w:r>' +
    '</w:p>'

$rng1 = $d.Paragraphs(1).Range
$rng1.InsertXML((New-PkgXml $para1))

# ---------------------------------------------------------------
# Paragraph 5 : "Quien es Jesus para ti?"
# ---------------------------------------------------------------
$para5 = '<w:p w14:paraId="374E7ABA" w14:textId="77777777" w:rsidR="001F22DF" w:rsidRPr="001F22DF" w:rsidRDefault="001F22DF" w:rsidP="001F22DF">' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&#191;Qui&#233;n es Jes&#250;s para t&#237;?</w:t></w:r>' +
    '</w:p>'

$rng5 = $d.Paragraphs(5).Range
$rng5.InsertXML((New-PkgXml $para5))

# ---------------------------------------------------------------
# Paragraph 21 : "Esta pregunta tambien nos la hace..."
# ---------------------------------------------------------------
$para21 = '<w:p w14:paraId="07CD849A" w14:textId="77777777" w:rsidR="001F22DF" w:rsidRPr="001F22DF" w:rsidRDefault="001F22DF" w:rsidP="001F22DF">' +
    '<w:r w:rsidRPr="001F22DF"><w:t xml:space="preserve">Esta pregunta tambi&#233;n nos la hace a nosotros hoy, me la hace a m&#237;, te la hace a t&#237;.  </w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">&#191;Qui&#233;n es Jes&#250;s para t&#237;? </w:t></w:r>' +
    '<w:r w:rsidRPr="001F22DF"><w:t>T&#243;mate unos segundos y responde mentalmente esta pregunta.  (dar unos segundos en silencio)</w:t></w:r>' +
    '</w:p>'

$rng21 = $d.Paragraphs(21).Range
$rng21.InsertXML((New-PkgXml $para21))

Write-Output "done"
